# LOQ4056.docx - "Requisitos" bullet list: move the "LOB1012 - Estatistica"
# requirement line from the first position to the last position, i.e.:
#   LOB1012, LOQ4095, LOQ4098   -->   LOQ4095, LOQ4098, LOB1012
#
# Each line in that bulleted paragraph lives in its own <w:r> (text run
# followed by a <w:br/>), so the edit below deletes the paragraph's run
# content and re-inserts the same three lines, reordered, as three
# separate runs (so the OOXML keeps one <w:r> per line, matching the
# original authoring style instead of collapsing everything into one run).

$d = $word.ActiveDocument

$line1 = "LOB1012 -  Estatística  (Requisito fraco)"
$line2 = "LOQ4095 -  Química Geral Experimental  (Requisito fraco)"
$line3 = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)"

# Locate the target paragraph robustly (don't hardcode a paragraph index):
# search for the distinguishing "LOB1012" requirement line, then expand the
# found hit to the whole enclosing paragraph (wdParagraph = 4).
$rng = $d.Content
$found = $rng.Find.Execute($line1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'LOB1012' requirement line to reorder"
}
$rng.Expand(4)

$paraStart = $rng.Start
$paraEnd = $rng.End

# Sanity-check the paragraph still has the exact three lines in the
# original order (text + manual line breaks, Chr(11), then the paragraph
# mark, Chr(13)) before touching anything.
$expected = $line1 + [char]11 + $line2 + [char]11 + $line3 + [char]11 + [char]13
if ($rng.Text -ne $expected) {
    throw "Requisitos paragraph did not match the expected layout; aborting"
}

# Remove the run content but keep the trailing paragraph mark (so the
# paragraph, its style/pPr and the rest of the document are untouched).
$contentRange = $d.Range($paraStart, $paraEnd - 1)
$contentRange.Delete()

# Re-insert the three lines in the new order as three distinct runs.
$newXml = '<?xml version="1.0" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:r><w:t>' + $line2 + '</w:t><w:br/></w:r>' + `
              '<w:r><w:t>' + $line3 + '</w:t><w:br/></w:r>' + `
              '<w:r><w:t>' + $line1 + '</w:t><w:br/></w:r>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$insertPoint = $d.Range($paraStart, $paraStart)
$insertPoint.InsertXML($newXml)

Write-Host ("Requisitos line reordered. New text: " + $rng.Text)
